$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3: enable review-option fields for this cycle (based on reviewCycle)
$ws.Range("B3").Value = "Senior manager feedback"
$ws.Range("C3").Value = "q1"
$ws.Range("D3").Value = 3
$ws.Range("F3").Value = 3
$ws.Range("G3").Value = 45321
$ws.Range("H3").Value = 3
$ws.Range("I3").Value = "Super senior manager feedback"
$ws.Range("K3").Value = "direct manager feedback"
$ws.Range("L3").Value = "Reviewed"
